# Update view-count / price figures across sheets as published on gh-pages
# (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 2789
$ws1.Range("F7").Value  = 2991
$ws1.Range("F8").Value  = 1907
$ws1.Range("F11").Value = 2572
$ws1.Range("F13").Value = 271
$ws1.Range("G17").Value = 61.92
$ws1.Range("F18").Value = 9546
$ws1.Range("F21").Value = 7536
$ws1.Range("F22").Value = 12073
$ws1.Range("F28").Value = 2717
$ws1.Range("F31").Value = 2708
$ws1.Range("F32").Value = 1057
$ws1.Range("F37").Value = 1118
$ws1.Range("F38").Value = 34
$ws1.Range("F41").Value = 576

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 198

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 2789
$ws4.Range("F11").Value = 2991
$ws4.Range("F13").Value = 1907
$ws4.Range("F16").Value = 2572
$ws4.Range("F19").Value = 271
$ws4.Range("G22").Value = 61.92
$ws4.Range("F23").Value = 9546
$ws4.Range("F25").Value = 7536
$ws4.Range("F26").Value = 12073
$ws4.Range("F33").Value = 2717
$ws4.Range("F46").Value = 576
